$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    "D8" = 5566400
    "E8" = 4812900
    "F8" = 5606500
    "G8" = 5009600
    "H8" = 3750100
    "I8" = 3624500
    "J8" = 3177300
    "D9" = 4483300
    "E9" = 3699300
    "F9" = 4068000
    "G9" = 3282300
    "H9" = 2525000
    "I9" = 2612500
    "J9" = 2529300
    "D10" = 1083100
    "E10" = 1113700
    "F10" = 1538500
    "G10" = 1727300
    "H10" = 1225100
    "I10" = 1012000
    "J10" = 648000
    "J14" = 1800
    "D17" = 4889000
    "E17" = 4350100
    "F17" = 4400900
    "G17" = 3578900
    "H17" = 2631700
    "I17" = 2733700
    "J17" = 2202000
    "D18" = 677500
    "E18" = 462800
    "F18" = 1205600
    "G18" = 1430700
    "H18" = 1118500
    "I18" = 890800
    "J18" = 975300
    "D20" = 137700
    "E20" = 125900
    "F20" = 11700
    "G20" = -43600
    "H20" = 99800
    "I20" = 625800
    "J20" = 68500
    "D21" = 1027900
    "E21" = 797700
    "F21" = 1426300
    "G21" = 1587600
    "H21" = 1424500
    "I21" = 1707800
    "J21" = 1240600
    "D22" = 393200
    "E22" = 494400
    "F22" = 355400
    "G22" = 238700
    "H22" = 178900
    "I22" = 208000
    "J22" = 232500
    "D23" = 422000
    "E23" = 94300
    "F23" = 861800
    "G23" = 1148400
    "H23" = 1039400
    "I23" = 1308600
    "J23" = 811200
    "D24" = 165200
    "E24" = 8500
    "F24" = 228800
    "G24" = 344200
    "H24" = 243600
    "I24" = 213500
    "J24" = 192000
    "D26" = 256800
    "E26" = 85800
    "F26" = 633000
    "G26" = 804200
    "H26" = 795800
    "I26" = 1095200
    "J26" = 619200
    "D27" = 256600
    "E27" = 85700
    "F27" = 632900
    "G27" = 804200
    "H27" = 795800
    "I27" = 1095200
    "J27" = 619200
    "D32" = -137700
    "E32" = -125900
    "F32" = -11700
    "G32" = 43600
    "H32" = -99800
    "I32" = -625800
    "J32" = -68500
    "D33" = 256600
    "E33" = 85700
    "F33" = 632900
    "G33" = 804200
    "H33" = 795800
    "I33" = 1095200
    "J33" = 619200
    "D35" = 256600
    "E35" = 85700
    "F35" = 632900
    "G35" = 804200
    "H35" = 795800
    "I35" = 1095200
    "J35" = 619200
    "D41" = 264100
    "E41" = 255100
    "F41" = 237100
    "G41" = 227400
    "H41" = 564500
    "I41" = 492000
    "J41" = 539400
    "D42" = 298600
    "E42" = 354200
    "F42" = 622200
    "G42" = 255100
    "H42" = 239800
    "I42" = 601800
    "J42" = 102900
    "D43" = 1315800
    "E43" = 1262800
    "F43" = 2235600
    "G43" = 881700
    "H43" = 850400
    "I43" = 1473100
    "J43" = 976000
    "D44" = 9800
    "E44" = 12700
    "F44" = 9600
    "G44" = 10200
    "H44" = 9700
    "I44" = 10600
    "J44" = 7900
    "D45" = 300500
    "E45" = 239400
    "F45" = 264700
    "G45" = 305900
    "H45" = 45300
    "I45" = 79000
    "J45" = 287900
    "D46" = 2188800
    "E46" = 2124200
    "F46" = 2404000
    "G46" = 1680400
    "H46" = 1709700
    "I46" = 2257100
    "J46" = 1478800
    "D47" = 4588500
    "E47" = 4391600
    "F47" = 4026800
    "G47" = 4685900
    "H47" = 3615700
    "I47" = 3712200
    "J47" = 3724200
    "D48" = 708200
    "E48" = 967900
    "F48" = 1010200
    "G48" = 1421300
    "H48" = 1491500
    "J48" = 3859500
    "D49" = 2860200
    "E49" = 2773900
    "F49" = 2634300
    "G49" = 866400
    "H49" = 513800
    "I49" = 961100
    "J49" = 2101100
    "D52" = 483700
    "E52" = 519600
    "F52" = 399600
    "G52" = 319400
    "H52" = 313000
    "I52" = 334300
    "J52" = 565800
    "D54" = 10829400
    "E54" = 10777200
    "F54" = 10474900
    "G54" = 8973300
    "H54" = 7643700
    "I54" = 8350300
    "J54" = 7594200
    "D57" = 600600
    "E57" = 497300
    "F57" = 487400
    "G57" = 411200
    "H57" = 273400
    "I57" = 334800
    "J57" = 216400
    "D58" = 607800
    "E58" = 1240100
    "F58" = 1615300
    "G58" = 1356400
    "H58" = 573700
    "I58" = 1657800
    "J58" = 2271100
    "D59" = 1012400
    "E59" = 1197500
    "F59" = 1249200
    "G59" = 827800
    "H59" = 671100
    "I59" = 1288600
    "J59" = 877400
    "D60" = 2220800
    "E60" = 2934900
    "F60" = 3351900
    "G60" = 2595400
    "H60" = 1518200
    "I60" = 3281200
    "J60" = 2076400
    "D61" = 3083500
    "E61" = 2651600
    "F61" = 2273100
    "G61" = 2107000
    "H61" = 1851000
    "I61" = 1012600
    "J61" = 1538300
    "D62" = 1851100
    "E62" = 1874600
    "F62" = 1520100
    "G62" = 1417800
    "H62" = 1034400
    "I62" = 1096700
    "J62" = 1356400
    "D66" = 7156500
    "E66" = 7462100
    "F66" = 7146200
    "G66" = 6081100
    "H66" = 4403500
    "I66" = 5389100
    "J66" = 4603300
    "D72" = 1468700
    "E72" = 1333100
    "F72" = 1195400
    "G72" = 665000
    "H72" = 984500
    "I72" = 2478200
    "J72" = 2702200
    "D76" = 3672900
    "E76" = 3315100
    "F76" = 3328800
    "G76" = 2892200
    "H76" = 3240200
    "I76" = 2961200
    "J76" = 2990900
    "D81" = 256600
    "E81" = 85700
    "F81" = 632900
    "G81" = 804200
    "H81" = 795800
    "I81" = 1095200
    "J81" = 619200
    "D83" = 217900
    "E83" = 213900
    "F83" = 214000
    "G83" = 205300
    "H83" = 211200
    "I83" = 195700
    "J83" = 201500
    "D89" = 148600
    "E89" = 311100
    "F89" = 771000
    "G89" = 957300
    "H89" = 901200
    "I89" = 622300
    "J89" = 939400
    "D91" = -21300
    "E91" = -30700
    "F91" = -32200
    "G91" = -31300
    "H91" = -17800
    "I91" = -27900
    "J91" = -24400
    "D94" = -99000
    "E94" = -157300
    "F94" = -824700
    "G94" = -1102100
    "H94" = 641600
    "I94" = -129500
    "J94" = -644800
    "D96" = -138400
    "E96" = -173000
    "F96" = -204200
    "G96" = -1004300
    "H96" = -1179400
    "I96" = -448000
    "J96" = -522000
    "D100" = -40600
    "E100" = -135700
    "F100" = 63300
    "G100" = -192200
    "H100" = -1470400
    "I100" = -540200
    "J100" = -453500
    "D102" = 9000
    "E102" = 18100
    "F102" = 9600
    "G102" = -337100
    "H102" = 72500
    "I102" = -47400
    "J102" = -159000
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
